# Updates the EC (Estado de Cuenta) worker data table on Hoja1.
# The data rows (B16:G29) are re-sorted from "grouped by worker" order
# (all periods for worker 1, then all periods for worker 2) to
# "interleaved by period" order (both workers' 2102, then both 2103, ...),
# and the "Valor Mora" amount (column G) is updated from 908526 to 877900
# for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Salario Basico, G=Valor Mora
$data = @(
  @(16, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2102", 36341, 877900),
  @(17, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2102", 36341, 877900),
  @(18, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2103", 36341, 877900),
  @(19, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2103", 36341, 877900),
  @(20, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2104", 36341, 877900),
  @(21, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2104", 36341, 877900),
  @(22, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2105", 36341, 877900),
  @(23, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2105", 36341, 877900),
  @(24, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2106", 36341, 877900),
  @(25, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2106", 36341, 877900),
  @(26, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2107", 36341, 877900),
  @(27, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2107", 36341, 877900),
  @(28, "CC", "1047420533", "ANDRES DE JESUS MARQUEZ LOBO", "2108", 33945, 877900),
  @(29, "CC", "9284806",    "JOSE SIMON RHENALS CASSIANI",  "2108", 35129, 877900)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
